$d = $word.ActiveDocument

$d.Content.Find.Execute("Presedenser er registrert.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Varsel: Presedenser er registrert.", 2)
